$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Paragraphs.Item(1).Range.Text = "2023-07-21 Friday"

# Update the division problems in the table, cell by cell
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "14÷9="
$t.Cell(1, 2).Range.Text = "37÷9="
$t.Cell(1, 3).Range.Text = "80÷5="
$t.Cell(1, 4).Range.Text = "91÷9="
$t.Cell(1, 5).Range.Text = "53÷9="
$t.Cell(5, 1).Range.Text = "64÷4="
$t.Cell(5, 2).Range.Text = "44÷3="
$t.Cell(5, 3).Range.Text = "63÷9="
$t.Cell(5, 4).Range.Text = "40÷4="
$t.Cell(5, 5).Range.Text = "79÷8="
$t.Cell(9, 1).Range.Text = "91÷7="
$t.Cell(9, 2).Range.Text = "70÷8="
$t.Cell(9, 3).Range.Text = "84÷2="
$t.Cell(9, 4).Range.Text = "90÷9="
$t.Cell(9, 5).Range.Text = "27÷9="
$t.Cell(13, 1).Range.Text = "39÷6="
$t.Cell(13, 2).Range.Text = "65÷2="
$t.Cell(13, 3).Range.Text = "39÷7="
$t.Cell(13, 4).Range.Text = "33÷9="
$t.Cell(13, 5).Range.Text = "94÷7="
$t.Cell(17, 1).Range.Text = "12÷8="
$t.Cell(17, 2).Range.Text = "47÷5="
$t.Cell(17, 3).Range.Text = "42÷8="
$t.Cell(17, 4).Range.Text = "42÷5="
$t.Cell(17, 5).Range.Text = "35÷8="
